$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3851.3635
$ws.Range("I64").Value = 2877.6667
$ws.Range("J64").Value = 4216.5
$ws.Range("K64").Value = 2877.6667
$ws.Range("L64").Value = 4216.5
$ws.Range("M64").Value = -2629.6667
$ws.Range("N64").Value = -4712.5

$ws.Range("H67").Value = 3851.3635
$ws.Range("I67").Value = 2877.6667
$ws.Range("J67").Value = 4216.5
$ws.Range("K67").Value = 2877.6667
$ws.Range("L67").Value = 4216.5
$ws.Range("M67").Value = -2019.6667
$ws.Range("N67").Value = -5932.5

$ws.Range("H131").Value = 6153.793
$ws.Range("I131").Value = 870
$ws.Range("J131").Value = 8166.6665
$ws.Range("K131").Value = 2610
$ws.Range("L131").Value = 24499.9995
$ws.Range("M131").Value = 2430
$ws.Range("N131").Value = -34579.99950000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 24244020
$ws.Range("I74").Value = 23810530
$ws.Range("J74").Value = 25644526
$ws.Range("K74").Value = 23810530
$ws.Range("L74").Value = 25644526
$ws.Range("M74").Value = -23809656
$ws.Range("N74").Value = -25646274

$ws.Range("H77").Value = 24244020
$ws.Range("I77").Value = 23810530
$ws.Range("J77").Value = 25644526
$ws.Range("K77").Value = 119052650
$ws.Range("L77").Value = 128222630
$ws.Range("M77").Value = -119048282
$ws.Range("N77").Value = -128231366

$ws.Range("H102").Value = 1917.4
$ws.Range("I102").Value = 1310.5714
$ws.Range("J102").Value = 3333.3333
$ws.Range("K102").Value = 1310.5714
$ws.Range("L102").Value = 3333.3333
$ws.Range("M102").Value = 311.4286
$ws.Range("N102").Value = -6577.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 526824.0600000001
$ws.Range("I107").Value = 556008.75
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 556008.75
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -554088.75
$ws.Range("N107").Value = -5340

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 83.333336
$ws.Range("I7").Value = 72.666664
$ws.Range("J7").Value = 94
$ws.Range("K7").Value = 72.666664
$ws.Range("L7").Value = 94
$ws.Range("M7").Value = 40.333336
$ws.Range("N7").Value = -320

$ws.Range("H31").Value = 1084361.6
$ws.Range("I31").Value = 1345071.9
$ws.Range("J31").Value = 6759.2666
$ws.Range("K31").Value = 1345071.9
$ws.Range("L31").Value = 6759.2666
$ws.Range("M31").Value = -1344776.9
$ws.Range("N31").Value = -7349.2666

$ws.Range("H34").Value = 1084361.6
$ws.Range("I34").Value = 1345071.9
$ws.Range("J34").Value = 6759.2666
$ws.Range("K34").Value = 1345071.9
$ws.Range("L34").Value = 6759.2666
$ws.Range("M34").Value = -1344869.9
$ws.Range("N34").Value = -7163.2666

$ws.Range("H62").Value = 2669.0952
$ws.Range("I62").Value = 2263
$ws.Range("J62").Value = 3684.3333
$ws.Range("K62").Value = 2263
$ws.Range("L62").Value = 3684.3333
$ws.Range("M62").Value = -1639
$ws.Range("N62").Value = -4932.3333

$ws.Range("H65").Value = 2669.0952
$ws.Range("I65").Value = 2263
$ws.Range("J65").Value = 3684.3333
$ws.Range("K65").Value = 11315
$ws.Range("L65").Value = 18421.6665
$ws.Range("M65").Value = -8195
$ws.Range("N65").Value = -24661.6665

$ws.Range("H132").Value = 1653.5
$ws.Range("I132").Value = 1227.85
$ws.Range("J132").Value = 2261.5715
$ws.Range("K132").Value = 3683.55
$ws.Range("L132").Value = 6784.7145
$ws.Range("M132").Value = -1153.55
$ws.Range("N132").Value = -11844.7145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 915
$ws.Range("I4").Value = 94.411766
$ws.Range("J4").Value = 7890
$ws.Range("K4").Value = 283.235298
$ws.Range("L4").Value = 23670
$ws.Range("M4").Value = -171.235298
$ws.Range("N4").Value = -23894

$ws.Range("H5").Value = 5209809
$ws.Range("I5").Value = 2175
$ws.Range("J5").Value = 6251336
$ws.Range("K5").Value = 6525
$ws.Range("L5").Value = 18754008
$ws.Range("M5").Value = -6413
$ws.Range("N5").Value = -18754232

$ws.Range("H132").Value = 2261.0833
$ws.Range("I132").Value = 935.6
$ws.Range("J132").Value = 3207.8572
$ws.Range("K132").Value = 8420.4
$ws.Range("L132").Value = 28870.7148
$ws.Range("M132").Value = -5890.4
$ws.Range("N132").Value = -33930.7148

$ws.Range("H135").Value = 5209809
$ws.Range("I135").Value = 2175
$ws.Range("J135").Value = 6251336
$ws.Range("K135").Value = 19575
$ws.Range("L135").Value = 56262024
$ws.Range("M135").Value = -17040
$ws.Range("N135").Value = -56267094

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8333.333000000001
$ws.Range("J5").Value = 8400
$ws.Range("L5").Value = 8400
$ws.Range("N5").Value = -8624

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 12024201
$ws.Range("I2").Value = 6000
$ws.Range("J2").Value = 20036334
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 20036334
$ws.Range("M2").Value = -5888
$ws.Range("N2").Value = -20036558

$ws.Range("H46").Value = 860.8484999999999
$ws.Range("I46").Value = 728.9167
$ws.Range("J46").Value = 936.2381
$ws.Range("K46").Value = 728.9167
$ws.Range("L46").Value = 936.2381
$ws.Range("M46").Value = -540.9167
$ws.Range("N46").Value = -1312.2381

$ws.Range("H68").Value = 3340.0833
$ws.Range("I68").Value = 1312.2
$ws.Range("J68").Value = 4788.5713
$ws.Range("K68").Value = 1312.2
$ws.Range("L68").Value = 4788.5713
$ws.Range("M68").Value = -563.2
$ws.Range("N68").Value = -6286.5713

$ws.Range("H71").Value = 3340.0833
$ws.Range("I71").Value = 1312.2
$ws.Range("J71").Value = 4788.5713
$ws.Range("K71").Value = 6561
$ws.Range("L71").Value = 23942.8565
$ws.Range("M71").Value = -2817
$ws.Range("N71").Value = -31430.8565

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9533.333000000001
$ws.Range("I2").Value = 9533.333000000001
$ws.Range("K2").Value = 9533.333000000001
$ws.Range("M2").Value = -9421.333000000001

$ws.Range("H122").Value = 4116678.8
$ws.Range("I122").Value = 6536904
$ws.Range("J122").Value = 2295.9
$ws.Range("K122").Value = 19610712
$ws.Range("L122").Value = 6887.700000000001
$ws.Range("M122").Value = -19608262
$ws.Range("N122").Value = -11787.7
